$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F336").Value = 102338
$ws.Range("F337").Value = 103523
$ws.Range("F338").Value = 227186
$ws.Range("G338").Value = 3190
$ws.Range("F339").Value = 661348
$ws.Range("F341").Value = 291420
$ws.Range("G341").Value = 3661
$ws.Range("F342").Value = 177851
$ws.Range("G342").Value = 3022
$ws.Range("F343").Value = 132946
$ws.Range("G343").Value = 2971
$ws.Range("F344").Value = 135007
$ws.Range("F345").Value = 291563
$ws.Range("F346").Value = 675375
$ws.Range("F347").Value = 343715
$ws.Range("F349").Value = 158926
$ws.Range("G349").Value = 2751
$ws.Range("F350").Value = 126800
$ws.Range("G350").Value = 2783
$ws.Range("F351").Value = 150714
$ws.Range("F352").Value = 307465
$ws.Range("F356").Value = 159817
$ws.Range("F358").Value = 158083
$ws.Range("F359").Value = 321242
$ws.Range("G359").Value = 3348
$ws.Range("F362").Value = 228296
$ws.Range("G362").Value = 3179
$ws.Range("F363").Value = 186952
$ws.Range("G363").Value = 2758
$ws.Range("F364").Value = 167498
$ws.Range("G364").Value = 2468
$ws.Range("F365").Value = 183859
$ws.Range("G365").Value = 2388
$ws.Range("F366").Value = 339997
$ws.Range("G366").Value = 2845
$ws.Range("F367").Value = 765513
$ws.Range("G367").Value = 3916
$ws.Range("F368").Value = 345789
$ws.Range("F369").Value = 234044
$ws.Range("G369").Value = 2605
$ws.Range("F370").Value = 182292
$ws.Range("G370").Value = 2042
$ws.Range("F371").Value = 159668
$ws.Range("G371").Value = 1956
$ws.Range("F372").Value = 178653
$ws.Range("G372").Value = 1860
$ws.Range("F373").Value = 348639
$ws.Range("G373").Value = 2375
$ws.Range("F374").Value = 771176
$ws.Range("G374").Value = 3420
$ws.Range("F376").Value = 221369
$ws.Range("G376").Value = 2232
$ws.Range("F377").Value = 176772
$ws.Range("G377").Value = 1825
$ws.Range("F378").Value = 157562
$ws.Range("G378").Value = 1554
$ws.Range("F379").Value = 179413
$ws.Range("F380").Value = 343669
$ws.Range("G380").Value = 1996
$ws.Range("F381").Value = 743254
$ws.Range("G381").Value = 2686
$ws.Range("F383").Value = 221218
$ws.Range("G383").Value = 1767
$ws.Range("F384").Value = 171797
$ws.Range("G384").Value = 1509
$ws.Range("F385").Value = 150613
$ws.Range("G385").Value = 1406
$ws.Range("F386").Value = 182429
$ws.Range("F387").Value = 351480
$ws.Range("G387").Value = 1677
$ws.Range("F388").Value = 718428
$ws.Range("G388").Value = 2157
$ws.Range("F389").Value = 350751
$ws.Range("G389").Value = 1297
$ws.Range("F390").Value = 218689
$ws.Range("G390").Value = 1506
$ws.Range("F391").Value = 175913
$ws.Range("G391").Value = 1214
$ws.Range("F392").Value = 214852
$ws.Range("G392").Value = 1191
$ws.Range("F393").Value = 290226
$ws.Range("G393").Value = 1171
$ws.Range("F394").Value = 158889
$ws.Range("G394").Value = 610
$ws.Range("F395").Value = 715814
$ws.Range("G395").Value = 1858
$ws.Range("F396").Value = 158921
$ws.Range("G396").Value = 531
$ws.Range("F397").Value = 102197
$ws.Range("G397").Value = 606
